$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- English-text QA pass over column C (translation fixes) ---
# Rows 1, 7, 14, 16, 23: plain text replacement + red font + wrap text (no vertical centering)
$r = $ws.Range("C1")
$r.VerticalAlignment = -4107
$r.WrapText = $true
$r.Font.Color = 255
$r.Value = "Operations"

$r = $ws.Range("C7")
$r.VerticalAlignment = -4107
$r.WrapText = $true
$r.Font.Color = 255
$r.Value = "Copied"

$r = $ws.Range("C14")
$r.VerticalAlignment = -4107
$r.WrapText = $true
$r.Font.Color = 255
$r.Value = "Modify description"

$r = $ws.Range("C16")
$r.VerticalAlignment = -4107
$r.WrapText = $true
$r.Font.Color = 255
$r.Value = "VPC"

$r = $ws.Range("C23")
$r.VerticalAlignment = -4107
$r.WrapText = $true
$r.Font.Color = 255
$r.Value = "Virtual Machine"

# Row 33: quote-prefixed text (forces text entry), default font, wrap + vertical center
$r = $ws.Range("C33")
$r.Value = "'Public IP deleted'"

# Row 36: quote-prefixed text, red font, wrap text (no vertical centering)
$r = $ws.Range("C36")
$r.VerticalAlignment = -4107
$r.WrapText = $true
$r.Font.Color = 255
$r.Value = "'Security Group Deleted'"
